$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "MLP(10,256)"
$ws.Range("H1").Value = "MLP(30,1024)"
$ws.Range("I1").Value = "xgboost"

# New data values for row 2 (20 news group)
$ws.Range("G2").Value = 0.7058
$ws.Range("H2").Value = 0.7408
$ws.Range("I2").Value = 0.6164

# New data values for row 3 (IMDB Reviews)
$ws.Range("G3").Value = 0.8757
$ws.Range("H3").Value = 0.8801
$ws.Range("I3").Value = 0.7398

# Set column width for new column F (target XML width ~16.35 chars)
$ws.Columns.Item(6).ColumnWidth = 15.45

# Set active selection cell (cosmetic change in the diff)
$ws.Range("H4").Select()
